# Auto-generated Excel COM-interop script
# Applies numeric updates to price/profit columns (H-N) in the Typhon_Profits workbook
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets, per the scheduled price refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H52").Value = 10000
$ws.Range("J52").Value = 10000
$ws.Range("L52").Value = 30000
$ws.Range("N52").Value = -30320

$ws.Range("H74").Value = 4052.9412
$ws.Range("I74").Value = 3300
$ws.Range("J74").Value = 4366.6665
$ws.Range("K74").Value = 3300
$ws.Range("L74").Value = 4366.6665
$ws.Range("M74").Value = -2364
$ws.Range("N74").Value = -6238.6665

$ws.Range("H77").Value = 4052.9412
$ws.Range("I77").Value = 3300
$ws.Range("J77").Value = 4366.6665
$ws.Range("K77").Value = 16500
$ws.Range("L77").Value = 21833.3325
$ws.Range("M77").Value = -11820
$ws.Range("N77").Value = -31193.3325

$ws.Range("H129").Value = 829.5
$ws.Range("J129").Value = 829.798
$ws.Range("L129").Value = 2489.394
$ws.Range("N129").Value = -12489.394

$ws.Range("H137").Value = 1989.9032
$ws.Range("I137").Value = 1857.2609
$ws.Range("K137").Value = 5571.7827
$ws.Range("M137").Value = -3021.7827

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2681.6726
$ws.Range("I32").Value = 2090.851
$ws.Range("K32").Value = 2090.851
$ws.Range("M32").Value = -1803.851

$ws.Range("H43").Value = 30375.666
$ws.Range("J43").Value = 30377
$ws.Range("L43").Value = 30377
$ws.Range("N43").Value = -31003

$ws.Range("H61").Value = 860174.3
$ws.Range("I61").Value = 1287879.4
$ws.Range("J61").Value = 4764.2856
$ws.Range("K61").Value = 1287879.4
$ws.Range("L61").Value = 4764.2856
$ws.Range("M61").Value = -1287667.4
$ws.Range("N61").Value = -5188.2856

$ws.Range("H74").Value = 2150.5
$ws.Range("I74").Value = 2117.2222
$ws.Range("J74").Value = 2450
$ws.Range("K74").Value = 2117.2222
$ws.Range("L74").Value = 2450
$ws.Range("M74").Value = -1243.2222
$ws.Range("N74").Value = -4198

$ws.Range("H77").Value = 2150.5
$ws.Range("I77").Value = 2117.2222
$ws.Range("J77").Value = 2450
$ws.Range("K77").Value = 10586.111
$ws.Range("L77").Value = 12250
$ws.Range("M77").Value = -6218.111000000001
$ws.Range("N77").Value = -20986

$ws.Range("H97").Value = 2425.4546
$ws.Range("I97").Value = 2088.889
$ws.Range("J97").Value = 3940
$ws.Range("K97").Value = 2088.889
$ws.Range("L97").Value = 3940
$ws.Range("M97").Value = -1592.889
$ws.Range("N97").Value = -4932

$ws.Range("H102").Value = 2642.6667
$ws.Range("I102").Value = 944.4545000000001
$ws.Range("K102").Value = 944.4545000000001
$ws.Range("M102").Value = 677.5454999999999

$ws.Range("H110").Value = 1751.8334
$ws.Range("I110").Value = 2070.1
$ws.Range("J110").Value = 1354
$ws.Range("K110").Value = 2070.1
$ws.Range("L110").Value = 1354
$ws.Range("M110").Value = -25.09999999999991
$ws.Range("N110").Value = -5444

$ws.Range("H132").Value = 15564.378
$ws.Range("I132").Value = 1711.579
$ws.Range("J132").Value = 30186.777
$ws.Range("K132").Value = 5134.737
$ws.Range("L132").Value = 90560.33099999999
$ws.Range("M132").Value = -2604.737
$ws.Range("N132").Value = -95620.33099999999

$ws.Range("H136").Value = 860174.3
$ws.Range("I136").Value = 1287879.4
$ws.Range("J136").Value = 4764.2856
$ws.Range("K136").Value = 3863638.2
$ws.Range("L136").Value = 14292.8568
$ws.Range("M136").Value = -3861088.2
$ws.Range("N136").Value = -19392.8568

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2233.3044
$ws.Range("I86").Value = 1834.3158
$ws.Range("J86").Value = 4128.5
$ws.Range("K86").Value = 1834.3158
$ws.Range("L86").Value = 4128.5
$ws.Range("M86").Value = -711.3158000000001
$ws.Range("N86").Value = -6374.5

$ws.Range("H89").Value = 2233.3044
$ws.Range("I89").Value = 1834.3158
$ws.Range("J89").Value = 4128.5
$ws.Range("K89").Value = 9171.579
$ws.Range("L89").Value = 20642.5
$ws.Range("M89").Value = -3555.579
$ws.Range("N89").Value = -31874.5

$ws.Range("H94").Value = 5361.4375
$ws.Range("I94").Value = 2823.2856
$ws.Range("K94").Value = 2823.2856
$ws.Range("M94").Value = -2372.2856

$ws.Range("H105").Value = 1908.2
$ws.Range("I105").Value = 1702.3
$ws.Range("J105").Value = 2320
$ws.Range("K105").Value = 1702.3
$ws.Range("L105").Value = 2320
$ws.Range("M105").Value = 44.70000000000005
$ws.Range("N105").Value = -5814

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7491.154
$ws.Range("I31").Value = 8350.953
$ws.Range("J31").Value = 3383.2222
$ws.Range("K31").Value = 8350.953
$ws.Range("L31").Value = 3383.2222
$ws.Range("M31").Value = -8055.953
$ws.Range("N31").Value = -3973.2222

$ws.Range("H34").Value = 7491.154
$ws.Range("I34").Value = 8350.953
$ws.Range("J34").Value = 3383.2222
$ws.Range("K34").Value = 8350.953
$ws.Range("L34").Value = 3383.2222
$ws.Range("M34").Value = -8148.953
$ws.Range("N34").Value = -3787.2222

$ws.Range("H57").Value = 12905.5
$ws.Range("J57").Value = 14707.333
$ws.Range("L57").Value = 14707.333
$ws.Range("N57").Value = -15827.333

$ws.Range("H58").Value = 26373.4
$ws.Range("I58").Value = 1529.4286
$ws.Range("J58").Value = 84342.664
$ws.Range("K58").Value = 1529.4286
$ws.Range("L58").Value = 84342.664
$ws.Range("M58").Value = -1326.4286
$ws.Range("N58").Value = -84748.664

$ws.Range("H134").Value = 1158.625
$ws.Range("I134").Value = 836.8461
$ws.Range("J134").Value = 1538.909
$ws.Range("K134").Value = 2510.5383
$ws.Range("L134").Value = 4616.727000000001
$ws.Range("M134").Value = 24.46169999999984
$ws.Range("N134").Value = -9686.727000000001

$ws.Range("H136").Value = 26373.4
$ws.Range("I136").Value = 1529.4286
$ws.Range("J136").Value = 84342.664
$ws.Range("K136").Value = 4588.2858
$ws.Range("L136").Value = 253027.992
$ws.Range("M136").Value = -2038.2858
$ws.Range("N136").Value = -258127.992

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1352.8667
$ws.Range("I5").Value = 1026.1818
$ws.Range("K5").Value = 3078.5454
$ws.Range("M5").Value = -2966.5454

$ws.Range("H131").Value = 741.65
$ws.Range("J131").Value = 755
$ws.Range("L131").Value = 2265
$ws.Range("N131").Value = -12345

$ws.Range("H132").Value = 925
$ws.Range("I132").Value = 925
$ws.Range("K132").Value = 8325
$ws.Range("M132").Value = -5795

$ws.Range("H135").Value = 1352.8667
$ws.Range("I135").Value = 1026.1818
$ws.Range("K135").Value = 9235.636200000001
$ws.Range("M135").Value = -6700.636200000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 4157
$ws.Range("I113").Value = 4000
$ws.Range("J113").Value = 4219.8
$ws.Range("K113").Value = 4000
$ws.Range("L113").Value = 4219.8
$ws.Range("M113").Value = -1830
$ws.Range("N113").Value = -8559.799999999999

$ws.Range("H126").Value = 3837.8918
$ws.Range("I126").Value = 2829.45
$ws.Range("J126").Value = 5024.294
$ws.Range("K126").Value = 8488.349999999999
$ws.Range("L126").Value = 15072.882
$ws.Range("M126").Value = -6018.349999999999
$ws.Range("N126").Value = -20012.882

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6877.6665
$ws.Range("I7").Value = 6875
$ws.Range("J7").Value = 6879.8
$ws.Range("K7").Value = 6875
$ws.Range("L7").Value = 6879.8
$ws.Range("M7").Value = -6763
$ws.Range("N7").Value = -7103.8

$ws.Range("H40").Value = 4623.357
$ws.Range("I40").Value = 3972.7
$ws.Range("J40").Value = 6250
$ws.Range("K40").Value = 3972.7
$ws.Range("L40").Value = 6250
$ws.Range("M40").Value = -3836.7
$ws.Range("N40").Value = -6522

$ws.Range("H46").Value = 2365.8333
$ws.Range("I46").Value = 2323.75
$ws.Range("J46").Value = 2450
$ws.Range("K46").Value = 2323.75
$ws.Range("L46").Value = 2450
$ws.Range("M46").Value = -2135.75
$ws.Range("N46").Value = -2826

$ws.Range("H126").Value = 6877.6665
$ws.Range("I126").Value = 6875
$ws.Range("J126").Value = 6879.8
$ws.Range("K126").Value = 20625
$ws.Range("L126").Value = 20639.4
$ws.Range("M126").Value = -18155
$ws.Range("N126").Value = -25579.4

$ws.Range("H132").Value = 2378.25
$ws.Range("I132").Value = 949.8570999999999
$ws.Range("J132").Value = 3489.2222
$ws.Range("K132").Value = 2849.5713
$ws.Range("L132").Value = 10467.6666
$ws.Range("M132").Value = -319.5712999999996
$ws.Range("N132").Value = -15527.6666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 2674815
$ws.Range("I107").Value = 327.1111
$ws.Range("J107").Value = 5683614
$ws.Range("K107").Value = 981.3333
$ws.Range("L107").Value = 17050842
$ws.Range("M107").Value = 938.6667
$ws.Range("N107").Value = -17054682

$ws.Range("H126").Value = 1713.8
$ws.Range("I126").Value = 967.6
$ws.Range("K126").Value = 2902.8
$ws.Range("M126").Value = -432.8000000000002

$ws.Range("H132").Value = 3297.0588
$ws.Range("I132").Value = 2923.0908
$ws.Range("J132").Value = 3982.6667
$ws.Range("K132").Value = 8769.2724
$ws.Range("L132").Value = 11948.0001
$ws.Range("M132").Value = -6239.2724
$ws.Range("N132").Value = -17008.0001
